# Log additional work on the JS101 time log: the 44497 ("1/8/2021") entry
# (row 61) gets extra minutes spent on Lesson 3 and the note is updated to
# reflect it. The weekly-total / grand-total formulas recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the note for the Finished-3-small-problems day to mention the extra
# time spent on Lesson 3.
$ws.Range("D61").Value = "Finished 3 small problems, 15 mins on Lesson 3"

# Add the 15 extra minutes (0.15 hr) onto the hours logged that day.
$ws.Range("C61").Value = 2.15

# Leave the selection where the author ended up after making the edit.
$ws.Range("A62").Select()
